$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write a value as literal TEXT (no Excel smart-parsing into
# number/date), by formatting a scratch cell as Text, writing the value,
# copying it (so the destination only receives the shared-string value,
# not the scratch cell's own formatting), then removing the scratch row.
function Set-TextValue {
    param($row, $col, [string]$text)

    $scratchRow = 500
    $scratchCell = $ws.Cells.Item($scratchRow, 1)
    $scratchCell.NumberFormat = "@"
    $scratchCell.Value = $text
    $ws.Range($scratchCell.Address()).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
    $ws.Rows($scratchRow).Delete()
}

# ---------------------------------------------------------------------
# Row 8  (Id 7)  - blank-style row, same pattern as existing row 4
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = 7
$ws.Range("H4:O4").Copy()
$ws.Range("H8:O8").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Row 9  (Id 8)  - blank-style row, same pattern as existing row 4
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = 8
$ws.Range("H4:O4").Copy()
$ws.Range("H9:O9").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Row 10 (Id 9)  - customer/dealer test data, blank stock/profit columns
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "shgda"
$ws.Cells.Item(10, 4).Value = "jlwdkas"
$ws.Cells.Item(10, 5).Value = "kjl"
Set-TextValue 10 6 "2018-09-15"
$ws.Range("H4:O4").Copy()
$ws.Range("H10:O10").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Row 11 (Id 10) - same customer/dealer data, blank stock/profit columns
# ---------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "shgda"
$ws.Cells.Item(11, 4).Value = "jlwdkas"
$ws.Cells.Item(11, 5).Value = "kjl"
Set-TextValue 11 6 "2018-09-15"
$ws.Range("H4:O4").Copy()
$ws.Range("H11:O11").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Row 12 (Id 11) - same customer/dealer data, WITH stock/profit filled in
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "shgda"
$ws.Cells.Item(12, 4).Value = "jlwdkas"
$ws.Cells.Item(12, 5).Value = "kjl"
Set-TextValue 12 6 "2018-09-15"
Set-TextValue 12 8 "1"
Set-TextValue 12 9 "034"
$ws.Cells.Item(12, 10).Value = 34
Set-TextValue 12 11 "031"
Set-TextValue 12 12 "354"
$ws.Cells.Item(12, 13).Value = 164.9
Set-TextValue 12 14 "131"
$ws.Cells.Item(12, 15).Value = 33.900000000000006

# ---------------------------------------------------------------------
# Row 13 (Id 12) - duplicate of row 12 data
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "shgda"
$ws.Cells.Item(13, 4).Value = "jlwdkas"
$ws.Cells.Item(13, 5).Value = "kjl"
$ws.Range("F12:O12").Copy()
$ws.Range("F13:O13").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Row 14 (Id 13) - duplicate of row 12 data
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "shgda"
$ws.Cells.Item(14, 4).Value = "jlwdkas"
$ws.Cells.Item(14, 5).Value = "kjl"
$ws.Range("F12:O12").Copy()
$ws.Range("F14:O14").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Row 15 (Id 14) - blank-style row, same pattern as existing row 4
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = 14
$ws.Range("H4:O4").Copy()
$ws.Range("H15:O15").PasteSpecial(-4163)
